# Append two new rows (3 and 4) of username/password test data to Sheet1,
# extending the used range from A1:C2 to A1:C4. Column C (Url) is left
# blank for these rows, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "balrajhatti32@gmail.com"
$ws.Range("B3").Value = "Opal@123"

$ws.Range("A4").Value = "sbasava022"
$ws.Range("B4").Value = "Opal@122"
